$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Output "Design.Name before: $($d.Name)"
try {
    $d.Name = "Office Theme"
    Write-Output "Set Design.Name ok, now: $($d.Name)"
} catch {
    Write-Output "ERROR: $_"
}
$sm = $d.SlideMaster
Write-Output "Master.Name: $($sm.Name)"
